# Update NATMI ligand-receptor score recomputation with new TPM-based values.
# The underlying per-cluster TPM expression values for ligand Col18a1 / receptor
# Itga5 changed, which changes every derived ligand/receptor average & total
# expression, specificity score, and edge weight/specificity in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 24.35712233333334
$ws.Range("H2").Value = 73.07136700000001
$ws.Range("I2").Value = 0.3750500562097488
$ws.Range("J2").Value = 0.3750500562097488
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 853.6416225850594
$ws.Range("R2").Value = 7682.774603265535
$ws.Range("S2").Value = 0.1434427657285723
$ws.Range("T2").Value = 0.1434427657285723
$ws.Range("G3").Value = 24.35712233333334
$ws.Range("H3").Value = 73.07136700000001
$ws.Range("I3").Value = 0.3750500562097488
$ws.Range("J3").Value = 0.3750500562097488
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 728.6011524229077
$ws.Range("R3").Value = 6557.41037180617
$ws.Range("S3").Value = 0.1224314298312616
$ws.Range("T3").Value = 0.1224314298312615
$ws.Range("G4").Value = 24.35712233333334
$ws.Range("H4").Value = 73.07136700000001
$ws.Range("I4").Value = 0.3750500562097488
$ws.Range("J4").Value = 0.3750500562097488
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 217.0640409783512
$ws.Range("R4").Value = 1953.576368805161
$ws.Range("S4").Value = 0.03647463473473302
$ws.Range("T4").Value = 0.03647463473473301
$ws.Range("G5").Value = 24.35712233333334
$ws.Range("H5").Value = 73.07136700000001
$ws.Range("I5").Value = 0.3750500562097488
$ws.Range("J5").Value = 0.3750500562097488
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 432.6519510338539
$ws.Range("R5").Value = 3893.867559304686
$ws.Range("S5").Value = 0.07270122591518191
$ws.Range("T5").Value = 0.0727012259151819
$ws.Range("I6").Value = 0.2805618708302703
$ws.Range("J6").Value = 0.2805618708302702
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 638.5795354130287
$ws.Range("R6").Value = 5747.215818717259
$ws.Range("S6").Value = 0.1073045318712589
$ws.Range("T6").Value = 0.1073045318712589
$ws.Range("I7").Value = 0.2805618708302703
$ws.Range("J7").Value = 0.2805618708302702
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("S7").Value = 0.09158668405231089
$ws.Range("T7").Value = 0.09158668405231087
$ws.Range("I8").Value = 0.2805618708302703
$ws.Range("J8").Value = 0.2805618708302702
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 162.3780410602207
$ws.Range("R8").Value = 1461.402369541986
$ws.Range("S8").Value = 0.02728540254718526
$ws.Range("T8").Value = 0.02728540254718525
$ws.Range("I9").Value = 0.2805618708302703
$ws.Range("J9").Value = 0.2805618708302702
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 323.6518400427567
$ws.Range("R9").Value = 2912.86656038481
$ws.Range("S9").Value = 0.05438525235951524
$ws.Range("T9").Value = 0.05438525235951523
$ws.Range("G10").Value = 22.31748066666667
$ws.Range("H10").Value = 66.952442
$ws.Range("I10").Value = 0.3436437303202491
$ws.Range("J10").Value = 0.343643730320249
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 782.1585057374399
$ws.Range("R10").Value = 7039.42655163696
$ws.Range("S10").Value = 0.1314310084381181
$ws.Range("T10").Value = 0.1314310084381181
$ws.Range("G11").Value = 22.31748066666667
$ws.Range("H11").Value = 66.952442
$ws.Range("I11").Value = 0.3436437303202491
$ws.Range("J11").Value = 0.343643730320249
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 667.5888025842994
$ws.Range("R11").Value = 6008.299223258694
$ws.Range("S11").Value = 0.1121791413147452
$ws.Range("T11").Value = 0.1121791413147452
$ws.Range("G12").Value = 22.31748066666667
$ws.Range("H12").Value = 66.952442
$ws.Range("I12").Value = 0.3436437303202491
$ws.Range("J12").Value = 0.343643730320249
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 198.8873099074318
$ws.Range("R12").Value = 1789.985789166886
$ws.Range("S12").Value = 0.03342028439879054
$ws.Range("T12").Value = 0.03342028439879052
$ws.Range("G13").Value = 22.31748066666667
$ws.Range("H13").Value = 66.952442
$ws.Range("I13").Value = 0.3436437303202491
$ws.Range("J13").Value = 0.343643730320249
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 396.4220986557011
$ws.Range("R13").Value = 3567.798887901311
$ws.Range("S13").Value = 0.06661329616859521
$ws.Range("T13").Value = 0.06661329616859518
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04834033333333334
$ws.Range("H14").Value = 0.145021
$ws.Range("I14").Value = 0.0007443426397318391
$ws.Range("J14").Value = 0.0007443426397318388
$ws.Range("M14").Value = 35.04689966666667
$ws.Range("N14").Value = 105.140699
$ws.Range("O14").Value = 0.3824629895491901
$ws.Range("P14").Value = 0.3824629895491901
$ws.Range("Q14").Value = 1.694178812186556
$ws.Range("R14").Value = 15.247609309679
$ws.Range("S14").Value = 0.000284683511240775
$ws.Range("T14").Value = 0.0002846835112407748
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04834033333333334
$ws.Range("H15").Value = 0.145021
$ws.Range("I15").Value = 0.0007443426397318391
$ws.Range("J15").Value = 0.0007443426397318388
$ws.Range("O15").Value = 0.3264402385872224
$ws.Range("P15").Value = 0.3264402385872223
$ws.Range("Q15").Value = 1.446017394549667
$ws.Range("R15").Value = 13.014156550947
$ws.Range("S15").Value = 0.0002429833889047044
$ws.Range("T15").Value = 0.0002429833889047043
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04834033333333334
$ws.Range("H16").Value = 0.145021
$ws.Range("I16").Value = 0.0007443426397318391
$ws.Range("J16").Value = 0.0007443426397318388
$ws.Range("M16").Value = 8.911727666666666
$ws.Range("N16").Value = 26.735183
$ws.Range("O16").Value = 0.09725271102035077
$ws.Range("P16").Value = 0.09725271102035075
$ws.Range("Q16").Value = 0.4307958859825556
$ws.Range("R16").Value = 3.877162973843
$ws.Range("S16").Value = 0.00007238933964196561
$ws.Range("T16").Value = 0.00007238933964196558
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.04834033333333334
$ws.Range("H17").Value = 0.145021
$ws.Range("I17").Value = 0.0007443426397318391
$ws.Range("J17").Value = 0.0007443426397318388
$ws.Range("M17").Value = 17.76285166666667
$ws.Range("N17").Value = 53.288555
$ws.Range("O17").Value = 0.1938440608432367
$ws.Range("P17").Value = 0.1938440608432367
$ws.Range("Q17").Value = 0.8586621705172223
$ws.Range("R17").Value = 7.727959534655001
$ws.Range("S17").Value = 0.0001442863999443941
$ws.Range("T17").Value = 0.000144286399944394

